$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New release: "Agosto.2021" column (BH) is added to the quarterly PIB
# staircase table. Every previously published quarter keeps the value last
# published for it (copied from column BG, "Mayo.2021") in the new column,
# except the most recently published quarter (row 74, 01-01-2021) which gets
# a revised figure, and a brand-new quarter (01-04-2021) is appended as a
# new row 75 with its first published value.
# ---------------------------------------------------------------------------

# 1) New header cell BH1 = "Agosto.2021", formatted like the rest of row 1
#    (bold, centered, bordered -> same style as BG1).
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("BH1").Value = "Agosto.2021"
$excel.CutCopyMode = 0

# 2) Carry forward the previously published values (rows 2 through 73) into
#    the new column BH, unchanged from column BG.
for ($r = 2; $r -le 73; $r++) {
    $ws.Cells.Item($r, 60).Value2 = $ws.Cells.Item($r, 59).Value2
}

# 3) Row 74 (01-01-2021) receives a revised value in the new column.
$ws.Cells.Item(74, 60).Value2 = 38076

# 4) New quarter 01-04-2021 is published for the first time -> new row 75.
#    Write the date-like label as a formula returning text, then collapse it
#    to a static value via paste-values; this stores it as plain text (not
#    an auto-parsed date serial) without leaving the cell's number format
#    changed, matching every other (unstyled) cell in column A.
$ws.Range("A75").Formula = "=""01-04-2021"""
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0

$ws.Cells.Item(75, 60).Value2 = 39677
